$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each row of this single-column table holds one metric value. Update the
# cells in place (by row index) so formatting (font/size) on the existing
# run is preserved, and collapse the three "raw sample dump" rows (44-46)
# down to the single summary value they should hold.

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "10842"

$t.Cell(6, 1).Range.Text  = "0.01822"
$t.Cell(7, 1).Range.Text  = "0.00371"
$t.Cell(8, 1).Range.Text  = "0.00017"
$t.Cell(9, 1).Range.Text  = "0.01740"
$t.Cell(10, 1).Range.Text = "0.01740"
$t.Cell(11, 1).Range.Text = "0.01822"
$t.Cell(12, 1).Range.Text = "2.00660"

$t.Cell(44, 1).Range.Text = "99.95"
$t.Cell(45, 1).Range.Text = "2.01"
$t.Cell(46, 1).Range.Text = "3900"
